$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Your palm is too close to the aircraft Please step farther away ."
$ws.Range("C2").Value = "Your palm is too close to the aircraft"
$ws.Range("D2").Value = "0-7"
$ws.Range("E2").Value = "Event"
$ws.Range("F2").Value = "Event"

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "Your palm is too close to the aircraft Please step farther away ."
$ws.Range("C3").Value = "Please step farther away"
$ws.Range("D3").Value = "8-11"
$ws.Range("E3").Value = "NonEvent"
$ws.Range("F3").Value = "NonEvent"

$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "RTK Base Station Data Not Received Please make sure base station antenna is connected to the correct port ."
$ws.Range("C4").Value = "RTK Base Station Data Not Received"
$ws.Range("D4").Value = "0-5"
$ws.Range("E4").Value = "Event"
$ws.Range("F4").Value = "Event"

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "RTK Base Station Data Not Received Please make sure base station antenna is connected to the correct port ."
$ws.Range("C5").Value = "Please make sure base station antenna is connected to the correct port"
$ws.Range("D5").Value = "6-17"
$ws.Range("E5").Value = "NonEvent"
$ws.Range("F5").Value = "NonEvent"

$ws.Range("A6").Value = 8
$ws.Range("B6").Value = "Warning: Critically low battery Please change the battery ."
$ws.Range("C6").Value = "Warning: Critically low battery"
$ws.Range("D6").Value = "0-3"
$ws.Range("E6").Value = "Event"
$ws.Range("F6").Value = "Event"

$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Warning: Critically low battery Please change the battery ."
$ws.Range("C7").Value = "Please change the battery"
$ws.Range("D7").Value = "4-7"
$ws.Range("E7").Value = "NonEvent"
$ws.Range("F7").Value = "NonEvent"

$ws.Range("A8").Value = 9
$ws.Range("B8").Value = "Vision sensor error Contact DJI Support for assistance ."
$ws.Range("C8").Value = "Contact DJI Support for assistance"
$ws.Range("D8").Value = "3-7"
$ws.Range("E8").Value = "NonEvent"
$ws.Range("F8").Value = "NonEvent"

$ws.Range("A9").Value = 22
$ws.Range("B9").Value = "Battery power low Check battery status and charge or warm up battery ."
$ws.Range("C9").Value = "Check battery status and charge or warm up battery"
$ws.Range("D9").Value = "3-11"
$ws.Range("E9").Value = "NonEvent"
$ws.Range("F9").Value = "NonEvent"

$ws.Range("A10").Value = 28
$ws.Range("B10").Value = "Your aircraft has entered a Warning Zone (Class D) Please fly with caution ."
$ws.Range("C10").Value = "Please fly with caution"
$ws.Range("D10").Value = "9-12"
$ws.Range("E10").Value = "NonEvent"
$ws.Range("F10").Value = "NonEvent"

$ws.Range("A11").Value = 33
$ws.Range("B11").Value = "Cannot track subject Subject too Small Get Closer and retry ."
$ws.Range("C11").Value = "Get Closer and retry"
$ws.Range("D11").Value = "6-9"
$ws.Range("E11").Value = "NonEvent"
$ws.Range("F11").Value = "NonEvent"

$ws.Range("A12").Value = 34
$ws.Range("B12").Value = "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant) Fly with caution ."
$ws.Range("C12").Value = "Fly with caution"
$ws.Range("D12").Value = "14-16"
$ws.Range("E12").Value = "NonEvent"
$ws.Range("F12").Value = "NonEvent"

$ws.Range("A13").Value = 36
$ws.Range("B13").Value = "Incompatible firmware version Go to Profile > Settings to update firmware ."
$ws.Range("C13").Value = "Go to Profile > Settings to update firmware"
$ws.Range("D13").Value = "3-10"
$ws.Range("E13").Value = "NonEvent"
$ws.Range("F13").Value = "NonEvent"

$ws.Range("A14").Value = 37
$ws.Range("B14").Value = "Failed to lower or pack up landing gear Take care of your gimbal when landing ."
$ws.Range("C14").Value = "Take care of your gimbal when landing"
$ws.Range("D14").Value = "8-14"
$ws.Range("E14").Value = "NonEvent"
$ws.Range("F14").Value = "NonEvent"

$ws.Range("A15").Value = 39
$ws.Range("B15").Value = "Battery overheating Stop flying and wait for battery temperature to return to normal ."
$ws.Range("C15").Value = "Stop flying and wait for battery temperature to return to normal"
$ws.Range("D15").Value = "2-12"
$ws.Range("E15").Value = "NonEvent"
$ws.Range("F15").Value = "NonEvent"

$ws.Range("A16").Value = 42
$ws.Range("B16").Value = "Aircraft in high interference environment Manually adjust flight route or return to home ."
$ws.Range("C16").Value = "Manually adjust flight route or return to home"
$ws.Range("D16").Value = "5-12"
$ws.Range("E16").Value = "NonEvent"
$ws.Range("F16").Value = "NonEvent"

$ws.Range("A17").Value = 45
$ws.Range("B17").Value = "Error: Course angle control error Please ensure the propellers are installed on the correct motors ."
$ws.Range("C17").Value = "Please ensure the propellers are installed on the correct motors"
$ws.Range("D17").Value = "5-14"
$ws.Range("E17").Value = "NonEvent"
$ws.Range("F17").Value = "NonEvent"

$ws.Range("A18").Value = 46
$ws.Range("B18").Value = "Weak signal Make sure the remote controller is facing aircraft and avoid blocking the signal ."
$ws.Range("C18").Value = "Weak signal"
$ws.Range("D18").Value = "0-1"
$ws.Range("E18").Value = "Event"
$ws.Range("F18").Value = "Event"
